$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds prices as plain text in the source data (e.g. "51.059.06",
# "12.40", "0.999"). Force Text format on every Column-D cell this script
# touches *before* writing its new value, so COM does not auto-coerce a
# numeric-looking string into a real number -- that coercion changes the cell
# type and can silently drop a significant trailing zero (e.g. "12.40" -> 12.4).
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '51.059.06'
$ws.Cells.Item(2, 5).Value = '  -0.96%  '

$ws.Cells.Item(3, 4).Value = '2.941.31'
$ws.Cells.Item(3, 5).Value = '  -1.26%  '

$ws.Cells.Item(4, 4).Value = '0.999'
$ws.Cells.Item(4, 5).Value = '  +0.00%  '

$ws.Cells.Item(5, 4).Value = '374.36'
$ws.Cells.Item(5, 5).Value = '  -1.58%  '

$ws.Cells.Item(6, 4).Value = '101.14'
$ws.Cells.Item(6, 5).Value = '  -2.59%  '

$ws.Cells.Item(7, 4).Value = '0.536'
$ws.Cells.Item(7, 5).Value = '  -1.39%  '

$ws.Cells.Item(8, 4).Value = '0.999'
$ws.Cells.Item(8, 5).Value = '  +0.00%  '

$ws.Cells.Item(9, 4).Value = '0.583'
$ws.Cells.Item(9, 5).Value = '  -1.95%  '

$ws.Cells.Item(10, 4).Value = '36.29'
$ws.Cells.Item(10, 5).Value = '  -2.71%  '

$ws.Cells.Item(11, 5).Value = '  -0.75%  '

$ws.Cells.Item(12, 4).Value = '0.0851'
$ws.Cells.Item(12, 5).Value = '  +0.49%  '

$ws.Cells.Item(13, 4).Value = '3.396.17'
$ws.Cells.Item(13, 5).Value = '  -1.46%  '

$ws.Cells.Item(14, 4).Value = '18.03'
$ws.Cells.Item(14, 5).Value = '  -2.26%  '

$ws.Cells.Item(15, 4).Value = '7.56'
$ws.Cells.Item(15, 5).Value = '  -0.28%  '

$ws.Cells.Item(16, 4).Value = '2.937.02'
$ws.Cells.Item(16, 5).Value = '  -1.33%  '

$ws.Cells.Item(17, 2).Value = 'Polygon'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(17, 4).Value = '0.994'
$ws.Cells.Item(17, 5).Value = '  +2.22%  '

$ws.Cells.Item(18, 2).Value = 'Uniswap'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Cells.Item(18, 4).Value = '10.69'
$ws.Cells.Item(18, 5).Value = '  +43.93%  '

$ws.Cells.Item(19, 4).Value = '50.967.45'
$ws.Cells.Item(19, 5).Value = '  -1.06%  '

$ws.Cells.Item(20, 4).Value = '3.11'
$ws.Cells.Item(20, 5).Value = '  -6.20%  '

$ws.Cells.Item(21, 4).Value = '12.40'
$ws.Cells.Item(21, 5).Value = '  -4.23%  '

$ws.Cells.Item(22, 4).Value = '0.0₃0956'
$ws.Cells.Item(22, 5).Value = '  -0.70%  '

$ws.Cells.Item(23, 4).Value = '265.63'
$ws.Cells.Item(23, 5).Value = '  +1.09%  '

$ws.Cells.Item(24, 4).Value = '68.61'
$ws.Cells.Item(24, 5).Value = '  -0.97%  '

$ws.Cells.Item(25, 4).Value = '3.12'
$ws.Cells.Item(25, 5).Value = '  +8.69%  '

$ws.Cells.Item(26, 4).Value = '8.08'
$ws.Cells.Item(26, 5).Value = '  -2.53%  '

$ws.Cells.Item(27, 4).Value = '7.65'
$ws.Cells.Item(27, 5).Value = '  -1.35%  '

$ws.Cells.Item(28, 5).Value = '  -0.01%  '

$ws.Cells.Item(29, 5).Value = '  -1.17%  '

$ws.Cells.Item(30, 5).Value = '  -4.37%  '

$ws.Cells.Item(31, 5).Value = '  -6.17%  '

$ws.Cells.Item(32, 4).Value = '10.01'
$ws.Cells.Item(32, 5).Value = '  +1.27%  '

$ws.Cells.Item(33, 4).Value = '50.67'
$ws.Cells.Item(33, 5).Value = '  -0.90%  '

$ws.Cells.Item(34, 5).Value = '  -1.24%  '

$ws.Cells.Item(35, 4).Value = '33.32'
$ws.Cells.Item(35, 5).Value = '  -4.97%  '

$ws.Cells.Item(36, 5).Value = '  -1.71%  '

$ws.Cells.Item(37, 5).Value = '  -0.21%  '

$ws.Cells.Item(38, 4).Value = '3.18'
$ws.Cells.Item(38, 5).Value = '  +4.67%  '

$ws.Cells.Item(39, 5).Value = '  -0.81%  '

$ws.Cells.Item(40, 4).Value = '16.22'
$ws.Cells.Item(40, 5).Value = '  -5.37%  '

$ws.Cells.Item(41, 2).Value = 'ARBITRUM'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Cells.Item(41, 4).Value = '1.79'
$ws.Cells.Item(41, 5).Value = '  -2.95%  '

$ws.Cells.Item(42, 2).Value = 'Stacks'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(42, 4).Value = '2.49'
$ws.Cells.Item(42, 5).Value = '  -3.71%  '

$ws.Cells.Item(43, 4).Value = '120.06'
$ws.Cells.Item(43, 5).Value = '  -4.93%  '

$ws.Cells.Item(44, 4).Value = '21.39'
$ws.Cells.Item(44, 5).Value = '  -1.26%  '

$ws.Cells.Item(45, 5).Value = '  -1.24%  '

$ws.Cells.Item(46, 5).Value = '  +2.37%  '

$ws.Cells.Item(47, 4).Value = '0.271'
$ws.Cells.Item(47, 5).Value = '  -4.02%  '

$ws.Cells.Item(48, 4).Value = '2.30'
$ws.Cells.Item(48, 5).Value = '  -3.13%  '

$ws.Cells.Item(49, 4).Value = '1.997.84'
$ws.Cells.Item(49, 5).Value = '  -1.92%  '

$ws.Cells.Item(50, 5).Value = '  -2.12%  '

$ws.Cells.Item(51, 5).Value = '  +1.95%  '
